$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Replace the "Railways" row (row 2) with a "Daycares" row.
$ws.Range("A2").Value = "Daycares"
$ws.Range("B2").Value = "All registered kindergartens in Norway, including public and private providers."
$ws.Range("C2").Value = "The Norwegian Directorate For Education And Training"
$ws.Range("D2").Value = "Annual"
$ws.Range("E2").Value = "High"
$ws.Range("F2").Value = "High"

# Row 2 no longer needs the taller 30pt row height used by the old
# (wrapped, multi-line) Railways text; restore it to the default/autofit
# height so no explicit row height is stored.
$ws.Rows.Item(2).AutoFit()

# Move the active selection.
$ws.Range("C17").Select()
